$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.135.49"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").Value = "2.451.85"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'568.53"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").Value = "'167.12"
$ws.Range("E6").Value = "  +4.46%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("E9").Value = "  +10.92%  "
$ws.Range("D10").Value = "2.451.75"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").Value = "'0.336"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").Value = "'4.71"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("E14").Value = "  +7.42%  "
$ws.Range("D15").Value = "70.104.27"
$ws.Range("E15").Value = "  +3.29%  "
$ws.Range("D16").Value = "2.904.01"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'24.23"
$ws.Range("E17").Value = "  +5.15%  "
$ws.Range("D18").Value = "2.454.45"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "'10.89"
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("D20").Value = "'7.15"
$ws.Range("E20").Value = "  +4.62%  "
$ws.Range("D21").Value = "'341.69"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("D23").Value = "'2.01"
$ws.Range("E23").Value = "  +8.27%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'66.41"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  +5.62%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.580.57"
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'8.52"
$ws.Range("E28").Value = "  +5.00%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "0.0₃0860"
$ws.Range("E30").Value = "  +6.23%  "
$ws.Range("E31").Value = "  +4.21%  "
$ws.Range("D32").Value = "'462.79"
$ws.Range("E32").Value = "  +10.17%  "
$ws.Range("E33").Value = "  +10.61%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("D36").Value = "'159.33"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  +6.98%  "
$ws.Range("D38").Value = "'19.10"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'18.25"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("D44").Value = "'38.10"
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("D46").Value = "'2.12"
$ws.Range("E46").Value = "  +5.60%  "
$ws.Range("D47").Value = "'134.62"
$ws.Range("E47").Value = "  +3.92%  "
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").Value = "'0.0726"
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("E51").Value = "  +1.70%  "
